$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

$startPara = $d.Paragraphs.Item($count - 2)
$endPara = $d.Paragraphs.Item($count)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
